# Generate Report for Handoff
# - Flip the localization status from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview + per-locale sheets).
# - Refresh the associated timestamps to reflect the new handoff generation.
# - Re-fit the now-shorter "Status" columns (Overview's zh-cn/de-de columns
#   and each locale sheet's Status column) to their new, narrower content.

$wb = $excel.ActiveWorkbook

# ColumnWidth (character units) set via COM is re-quantized by the host to
# the nearest 1/6 of a character (pixel-snap at MDW=6) before being written
# out as the OOXML <col width="..."> attribute, where width = ColumnWidth +
# 5/6. Back-solve the ColumnWidth to feed in so the stored width lands as
# close as possible to the desired value.
function Get-ColumnWidthInput($targetStoredWidth) {
    $mdw = 6.0
    $desired = $targetStoredWidth - (5.0 / $mdw)
    $steps = [Math]::Round($desired * $mdw)
    return $steps / $mdw
}

$newStatus = "Ready for handoff"
$newStatusColWidth = Get-ColumnWidthInput 17.2159881591797

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-23 13:00:47"
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-23 13:00:37"
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
